$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the date-looking label as text (not auto-converted to a date serial)
# by entering it as a formula returning a string literal, then flattening
# the formula to its cached value via copy / paste-special values.
$ws.Cells.Item(52, 1).Formula = "=""01-07-2021"""
$ws.Cells.Item(52, 1).Copy() | Out-Null
$ws.Cells.Item(52, 1).PasteSpecial(-4163) | Out-Null

$ws.Cells.Item(52, 2).Value = 118417
